# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Cells in column D that look like plain decimals get NumberFormat "@" (Text)
# set first so Excel's automatic number inference doesn't coerce strings such
# as "579.84" into floating-point numbers (and mangle trailing zeros); values
# that already contain multiple dots (e.g. "67.135.59") or other non-numeric
# characters are left alone since Excel already stores them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.135.59'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '3.120.10'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.84'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.34'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.45'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.07'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").Value = '3.637.16'
$ws.Range("D16").Value = '67.111.78'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.12'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.121.36'
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.71'
$ws.Range("E19").Value = '  +3.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '491.80'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("E21").Value = '  +5.53%  '
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.49'
$ws.Range("E26").Value = '  +4.33%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.94'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.67'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.57'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").Value = '0.0₃0948'
$ws.Range("E33").Value = '  -5.44%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '47.35'
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("E40").Value = '  +1.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.51'
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("D42").Value = '2.822.34'
$ws.Range("E42").Value = '  -0.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '385.32'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("E44").Value = '  -6.63%  '
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.60'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.88'
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("E51").Value = '  -1.50%  '
